# Adding code for Cases functionality
$wb = $excel.ActiveWorkbook

# --- Update selection on the "Companies" sheet (sheet2) ---
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Activate()
$wsCompanies.Range("E26").Select() | Out-Null

# --- Update selection on the "Deals" sheet (sheet3): select entire row 1 ---
$wsDeals = $wb.Worksheets.Item("Deals")
$wsDeals.Activate()
$wsDeals.Rows.Item(1).Select() | Out-Null

# --- Add the new "Cases" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCases = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCases.Name = "Cases"

# Header row (column by column so the shared-string table is built in the
# same order the columns are populated)
$headers = @("title","status","identifier","type","priority","contact","tags","description")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $headerCell = $wsCases.Cells.Item(1, $c + 1)
    $headerCell.Value = $headers[$c]
    $headerCell.Interior.Color = 65535
}

# Data rows, populated column by column
$columns = @(
    @("CaseTitle1", "CaseTitle2"),
    @("Awaiting input", "Enquiring"),
    @("aaaa", "bbbb"),
    @("Business Support", "Complaint"),
    @("High", "Low"),
    @("aaaaa", "zzzx"),
    @("Case 1 tag", "Case 2 tag"),
    @("Case 1 description", "Case 2 description")
)

for ($c = 0; $c -lt $columns.Count; $c++) {
    $col = $columns[$c]
    for ($r = 0; $r -lt $col.Count; $r++) {
        $wsCases.Cells.Item($r + 2, $c + 1).Value = $col[$r]
    }
}

$headerRange = $wsCases.Range($wsCases.Cells.Item(1, 1), $wsCases.Cells.Item(1, $headers.Count))
$headerRange.EntireColumn.AutoFit() | Out-Null
$wsCases.Range("E3").Select() | Out-Null
